$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 7583.1665
$ws.Range("I74").Value = 7299.8
$ws.Range("K74").Value = 7299.8
$ws.Range("M74").Value = -6363.8

$ws.Range("H77").Value = 7583.1665
$ws.Range("I77").Value = 7299.8
$ws.Range("K77").Value = 36499
$ws.Range("M77").Value = -31819

$ws.Range("H106").Value = 648
$ws.Range("I106").Value = 648
$ws.Range("K106").Value = 648
$ws.Range("M106").Value = -17

$ws.Range("H113").Value = 4497.222
$ws.Range("I113").Value = 4496.4
$ws.Range("J113").Value = 4498.25
$ws.Range("K113").Value = 4496.4
$ws.Range("L113").Value = 4498.25
$ws.Range("M113").Value = -1242.4
$ws.Range("N113").Value = -11006.25

$ws.Range("H137").Value = 1912.55
$ws.Range("I137").Value = 1031.4546
$ws.Range("K137").Value = 3094.3638
$ws.Range("M137").Value = -544.3638000000001


# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3780.913
$ws.Range("I74").Value = 3565.0476
$ws.Range("J74").Value = 6047.5
$ws.Range("K74").Value = 3565.0476
$ws.Range("L74").Value = 6047.5
$ws.Range("M74").Value = -2691.0476
$ws.Range("N74").Value = -7795.5

$ws.Range("H76").Value = 29988
$ws.Range("J76").Value = 29988
$ws.Range("L76").Value = 29988
$ws.Range("N76").Value = -30664

$ws.Range("H77").Value = 3780.913
$ws.Range("I77").Value = 3565.0476
$ws.Range("J77").Value = 6047.5
$ws.Range("K77").Value = 17825.238
$ws.Range("L77").Value = 30237.5
$ws.Range("M77").Value = -13457.238
$ws.Range("N77").Value = -38973.5

$ws.Range("H79").Value = 29988
$ws.Range("J79").Value = 29988
$ws.Range("L79").Value = 29988
$ws.Range("N79").Value = -32328

$ws.Range("H122").Value = 2248.8462
$ws.Range("I122").Value = 1798.5
$ws.Range("K122").Value = 5395.5
$ws.Range("M122").Value = -2945.5

$ws.Range("H132").Value = 1747.6888
$ws.Range("I132").Value = 1760.619
$ws.Range("K132").Value = 5281.857
$ws.Range("M132").Value = -2751.857


# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1324.0741
$ws.Range("I134").Value = 1155.8077
$ws.Range("K134").Value = 3467.4231
$ws.Range("M134").Value = -932.4231


# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H104").Value = 40095
$ws.Range("J104").Value = 40095
$ws.Range("L104").Value = 40095
$ws.Range("N104").Value = -45337

$ws.Range("H105").Value = 1434
$ws.Range("I105").Value = 1270.7142
$ws.Range("K105").Value = 1270.7142
$ws.Range("M105").Value = 476.2858000000001

$ws.Range("H109").Value = 46081.43
$ws.Range("J109").Value = 46081.43
$ws.Range("L109").Value = 46081.43
$ws.Range("N109").Value = -48161.43

$ws.Range("H132").Value = 3330
$ws.Range("I132").Value = 3173.25
$ws.Range("J132").Value = 4897.5
$ws.Range("K132").Value = 9519.75
$ws.Range("L132").Value = 14692.5
$ws.Range("M132").Value = -6989.75
$ws.Range("N132").Value = -19752.5

$ws.Range("H134").Value = 1524.5625
$ws.Range("I134").Value = 1602.7693
$ws.Range("K134").Value = 4808.3079
$ws.Range("M134").Value = -2273.3079

$ws.Range("H141").Value = 82200
$ws.Range("J141").Value = 82200
$ws.Range("L141").Value = 82200
$ws.Range("N141").Value = -92560


# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 869.44446
$ws.Range("I122").Value = 865.5
$ws.Range("K122").Value = 7789.5
$ws.Range("M122").Value = -5339.5


# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3038.2307
$ws.Range("I126").Value = 3038.2307
$ws.Range("K126").Value = 9114.6921
$ws.Range("M126").Value = -6644.6921

$ws.Range("H132").Value = 36454.234
$ws.Range("I132").Value = 47526.363
$ws.Range("J132").Value = 6005.875
$ws.Range("K132").Value = 142579.089
$ws.Range("L132").Value = 18017.625
$ws.Range("M132").Value = -140049.089
$ws.Range("N132").Value = -23077.625


# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4400.4
$ws.Range("J16").Value = 2000
$ws.Range("L16").Value = 2000
$ws.Range("N16").Value = -2340

$ws.Range("H29").Value = 3000
$ws.Range("I29").Value = 3000
$ws.Range("K29").Value = 3000
$ws.Range("M29").Value = -2705

$ws.Range("H40").Value = 5666.3335
$ws.Range("I40").Value = 4502
$ws.Range("K40").Value = 4502
$ws.Range("M40").Value = -4366

$ws.Range("H53").Value = 7311.75
$ws.Range("I53").Value = 7311.75
$ws.Range("K53").Value = 7311.75
$ws.Range("M53").Value = -6793.75

$ws.Range("H93").Value = 1031.48
$ws.Range("I93").Value = 970.5714
$ws.Range("J93").Value = 1351.25
$ws.Range("K93").Value = 970.5714
$ws.Range("L93").Value = 1351.25
$ws.Range("M93").Value = 277.4286
$ws.Range("N93").Value = -3847.25


# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 1562.5
$ws.Range("I2").Value = 2038.3334
$ws.Range("J2").Value = 135
$ws.Range("K2").Value = 2038.3334
$ws.Range("L2").Value = 135
$ws.Range("M2").Value = -1926.3334
$ws.Range("N2").Value = -359

$ws.Range("H9").Value = 5000
$ws.Range("I9").Value = 5000
$ws.Range("K9").Value = 5000
$ws.Range("M9").Value = -4860

$ws.Range("H13").Value = 979
$ws.Range("I13").Value = 979
$ws.Range("K13").Value = 979
$ws.Range("M13").Value = -839

$ws.Range("H62").Value = 10054.889
$ws.Range("J62").Value = 10928.429
$ws.Range("L62").Value = 10928.429
$ws.Range("N62").Value = -12176.429

$ws.Range("H65").Value = 10054.889
$ws.Range("J65").Value = 10928.429
$ws.Range("L65").Value = 54642.145
$ws.Range("N65").Value = -60882.145

$ws.Range("H122").Value = 6070
$ws.Range("I122").Value = 2500
$ws.Range("K122").Value = 7500
$ws.Range("M122").Value = -5050

$ws.Range("H126").Value = 3036.84
$ws.Range("I126").Value = 1338.8125
$ws.Range("K126").Value = 4016.4375
$ws.Range("M126").Value = -1546.4375

$ws.Range("H135").Value = 70000
$ws.Range("J135").Value = 70000
$ws.Range("L135").Value = 70000
$ws.Range("N135").Value = -80140

$ws.Range("H136").Value = 3765.5
$ws.Range("J136").Value = 5433
$ws.Range("L136").Value = 16299
$ws.Range("N136").Value = -21399

